# Solved the "Word Break" LeetCode problem -> add a new tracking row
# (date, question title, URL) right after the existing "House Robber" row,
# together with its URL hyperlink, mirroring the existing rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet currently ends at row 19 (House Robber). Row 19's date cell
# (A19) happens to use a slightly different - but visually identical -
# date style than the rest of the date column. When the new row 20 is
# appended, that "odd one out" style moves down to the new row, and A19
# falls back in line with the rest of the column (same style as A18).
# ------------------------------------------------------------------

# Carry A19's current (odd) date formatting down onto the new A20 cell
# before we touch A19 itself.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = 46080

# Re-format A19 so that it matches the rest of the date column (A18's style).
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)

# New question/title text.
$ws.Range("B20").Value = "Word Break"

# New URL cell, plus matching hyperlink relationship.
$ws.Range("C20").Value = "https://leetcode.com/problems/word-break/description/"
$ws.Hyperlinks.Add($ws.Range("C20"), "https://leetcode.com/problems/word-break/description/")

# Re-apply the same visual style the other URL cells use (Excel's
# Hyperlinks.Add forces its own built-in "Hyperlink" style onto the cell;
# put it back in line with the rest of column C).
$ws.Range("C19").Copy()
$ws.Range("C20").PasteSpecial(-4122)

# Drop the now-unused built-in "Hyperlink" named style that Hyperlinks.Add
# created, since the workbook does not otherwise use it.
$wb.Styles.Item("Hyperlink").Delete()

Write-Host "Added Word Break row (row 20)"
